$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data table: each entry is Row, then values for columns B,C,D,E,F,I,J,K,L,M,N
$data = @(
    @(2, 1.02, 1.031673357105697, 1.036299454614937, 1.041561461487249, 1.053401397699271, 1.037183268039884, 1.036807787208039, 1.039093592904857, 1.044340627457229, 1.056147467704779, 1.016257579537647),
    @(3, 1.02, 1.032448907208464, 1.036879201612525, 1.04228647119365, 1.054289464434955, 1.037348341686292, 1.037226196426974, 1.039483522881351, 1.04487651486592, 1.056848379585562, 1.016397247667257),
    @(4, 1.02, 1.032951423219397, 1.037254877379135, 1.042756638167133, 1.054865440474069, 1.037454323471892, 1.037496932851472, 1.039735672757486, 1.04522363525291, 1.057302614283762, 1.016487592473831),
    @(5, 1.02, 1.033162842675664, 1.037412939052577, 1.042954542719373, 1.055107898905812, 1.037498678312138, 1.037610748603765, 1.039841636735397, 1.045369650622085, 1.057493739813561, 1.016525565827425),
    @(6, 1.02, 1.033198350362783, 1.037439485730284, 1.042987786154138, 1.055148627336456, 1.037506113941029, 1.037629858624078, 1.039859426179746, 1.045394172241549, 1.057525840269709, 1.016531941261182),
    @(7, 1.02, 1.032954247583355, 1.037256988907989, 1.04275928161189, 1.054868678968188, 1.037454916929845, 1.03749845367197, 1.039737088811497, 1.045225585981569, 1.057305167465338, 1.0164880999061),
    @(8, 1.02, 1.031935315317292, 1.036495269754603, 1.041806265975334, 1.053701245819862, 1.037239227254863, 1.036949190269691, 1.039225404091185, 1.044521656544314, 1.056384198579679, 1.01630478692887),
    @(9, 1.02, 1.030145140840554, 1.035157246633561, 1.040134960933275, 1.051654415480588, 1.036852816069308, 1.035981357099771, 1.038322576034159, 1.043284108676229, 1.054766758169876, 1.015981558311606),
    @(10, 1.02, 1.02895537447267, 1.034268182088157, 1.039026273612824, 1.050296936394121, 1.036590991587837, 1.035336238174656, 1.037719979676497, 1.042461094399615, 1.053692219525757, 1.015765960743085),
    @(11, 1.02, 1.028441089431579, 1.033883930844346, 1.038547532890268, 1.049710837439343, 1.03647662790605, 1.035056935310695, 1.037458896787798, 1.042105217938621, 1.053227845907692, 1.015672583987217),
    @(12, 1.02, 1.028250196915012, 1.033741312854924, 1.038369908839657, 1.049493391522262, 1.036433999864268, 1.034953196799804, 1.037361897081946, 1.041973105378518, 1.053055495248033, 1.01563789693627),
    @(13, 1.02, 1.028291137830181, 1.033771899887782, 1.038408000647571, 1.049540022730192, 1.036443150421759, 1.034975448724244, 1.037382704817405, 1.042001440515808, 1.053092458757878, 1.015645337543972),
    @(14, 1.02, 1.028425307414472, 1.033872139741133, 1.038532846307186, 1.04969285803005, 1.036473107279299, 1.035048360102434, 1.03745087920076, 1.042094295917801, 1.05321359651479, 1.015669716798439),
    @(15, 1.02, 1.028507991780913, 1.033933915460816, 1.038609794608896, 1.049787059065892, 1.036491545062121, 1.035093284177019, 1.037492880827107, 1.042151517286926, 1.053288251876861, 1.015684737309436),
    @(16, 1.02, 1.028989524824039, 1.034293698902547, 1.039058074214251, 1.050335869820021, 1.036598560693536, 1.03535477549651, 1.037737303743579, 1.042484723292846, 1.053723057777167, 1.015772157436763),
    @(17, 1.02, 1.029291817830305, 1.034519575559051, 1.03933962519305, 1.050680580552577, 1.036665423727641, 1.035518813165799, 1.037890583315785, 1.04269386801855, 1.053996044626976, 1.015826988330607),
    @(18, 1.02, 1.029468226130905, 1.034651394890621, 1.039503977132364, 1.050881808221604, 1.036704328044584, 1.035614497130338, 1.037979973681589, 1.042815906172345, 1.054155360790335, 1.015858968148636),
    @(19, 1.02, 1.029528391316139, 1.034696353571963, 1.039560038548772, 1.05095044936249, 1.036717577145952, 1.035647123454578, 1.038010450898508, 1.042857526050842, 1.054209698278915, 1.01586987206695),
    @(20, 1.02, 1.029259375752508, 1.034495333949652, 1.039309404181634, 1.050643579402283, 1.03665825985632, 1.035501213097842, 1.037874139411731, 1.04267142383677, 1.053966746652177, 1.015821105709905),
    @(21, 1.02, 1.028385794050808, 1.033842618549899, 1.038496076791586, 1.049647844722636, 1.036464289818574, 1.035026889314779, 1.037430804137719, 1.0420669502112, 1.053177920629835, 1.01566253778817),
    @(22, 1.02, 1.02783732546228, 1.03343286802887, 1.037985872800307, 1.04902327671453, 1.036341475637709, 1.034728704328425, 1.037151935634756, 1.041687332996563, 1.052682755934538, 1.015562824073427),
    @(23, 1.02, 1.028128003819002, 1.033650023458283, 1.038256230275736, 1.049354230086077, 1.036406662790281, 1.034886773486529, 1.037299780601179, 1.041888533098247, 1.052945175512437, 1.015615685535865),
    @(24, 1.02, 1.029274034665828, 1.034506287475609, 1.039323059360566, 1.050660298123789, 1.036661497197385, 1.035509165800118, 1.037881569750727, 1.04268156524963, 1.053979984875594, 1.015823763819583),
    @(25, 1.02, 1.03060730272643, 1.035502645711533, 1.040566069664905, 1.052182332426559, 1.036953459172250, 1.036231552614377, 1.038556109788468, 1.043603695494985, 1.055184251350117, 1.016065142403652)
)

$colIndexes = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)  # B, C, D, E, F, I, J, K, L, M, N

foreach ($entry in $data) {
    $rowNum = $entry[0]
    for ($i = 0; $i -lt $colIndexes.Length; $i++) {
        $colNum = $colIndexes[$i]
        $value = $entry[$i + 1]
        $ws.Cells.Item($rowNum, $colNum).Value = $value
    }
}

Write-Host "Done updating vm_pu values"